$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Altman's Z)
$ws.Range("B2").Value = 0.7415

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0.9118"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.7415"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = 0.815

# Row 3 (Financial Variables and Sector)
$ws.Range("B3").Value = 0.9436

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0.9120"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.9436"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = 0.9276
